$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'41.756.98"
$ws.Range("E2").Value = "  +2.72%  "
$ws.Range("D3").Value = "'2.266.08"
$ws.Range("E3").Value = "  +1.47%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'304.09"
$ws.Range("E5").Value = "  +0.68%  "
$ws.Range("D6").Value = "'92.13"
$ws.Range("E6").Value = "  +2.49%  "
$ws.Range("E7").Value = "  +2.38%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "'0.481"
$ws.Range("E9").Value = "  +0.78%  "
$ws.Range("D10").Value = "'32.32"
$ws.Range("E10").Value = "  +2.68%  "
$ws.Range("D11").Value = "'53.20"
$ws.Range("E11").Value = "  +0.64%  "
$ws.Range("D12").Value = "'0.0797"
$ws.Range("E12").Value = "  +1.03%  "
$ws.Range("E13").Value = "  +0.44%  "
$ws.Range("E14").Value = "  +1.74%  "
$ws.Range("D15").Value = "'2.621.05"
$ws.Range("E15").Value = "  +2.00%  "
$ws.Range("D16").Value = "'14.19"
$ws.Range("E16").Value = "  +1.33%  "
$ws.Range("D17").Value = "'2.269.38"
$ws.Range("E17").Value = "  +1.20%  "
$ws.Range("D18").Value = "'0.764"
$ws.Range("E18").Value = "  +2.60%  "
$ws.Range("D19").Value = "'41.661.85"
$ws.Range("E19").Value = "  +2.89%  "
$ws.Range("D20").Value = "'12.65"
$ws.Range("D21").Value = "'0.0₃0903"
$ws.Range("E21").Value = "  +0.88%  "
$ws.Range("E22").Value = "  +1.43%  "
$ws.Range("D23").Value = "'66.83"
$ws.Range("E23").Value = "  +1.21%  "
$ws.Range("D24").Value = "'239.87"
$ws.Range("E24").Value = "  +1.07%  "
$ws.Range("D25").Value = "'2.59"
$ws.Range("E25").Value = "  +1.86%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("E27").Value = "  +4.01%  "
$ws.Range("D28").Value = "'24.03"
$ws.Range("E28").Value = "  +0.92%  "
$ws.Range("E29").Value = "  +0.87%  "
$ws.Range("E30").Value = "  -4.39%  "
$ws.Range("D31").Value = "'160.60"
$ws.Range("E31").Value = "  +2.20%  "
$ws.Range("D32").Value = "'34.49"
$ws.Range("E32").Value = "  +4.94%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "'5.23"
$ws.Range("E33").Value = "  +3.61%  "
$ws.Range("B34").Value = "FirstDigitalUSD"
$ws.Range("C34").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D34").Value = "'0.999"
$ws.Range("E34").Value = "  -0.11%  "
$ws.Range("D35").Value = "'0.0744"
$ws.Range("E35").Value = "  +2.79%  "
$ws.Range("E36").Value = "  +0.11%  "
$ws.Range("B37").Value = "Celestia"
$ws.Range("C37").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D37").Value = "'16.90"
$ws.Range("E37").Value = "  +6.36%  "
$ws.Range("B38").Value = "WEMIXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D38").Value = "'2.38"
$ws.Range("E38").Value = "  +1.87%  "
$ws.Range("E40").Value = "  +0.66%  "
$ws.Range("D41").Value = "'1.80"
$ws.Range("E41").Value = "  +1.64%  "
$ws.Range("D42").Value = "'3.95"
$ws.Range("E42").Value = "  +2.16%  "
$ws.Range("D43").Value = "'2.037.26"
$ws.Range("E43").Value = "  -2.22%  "
$ws.Range("D44").Value = "'19.29"
$ws.Range("E44").Value = "  -1.85%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "'0.0278"
$ws.Range("E45").Value = "  +1.80%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").Value = "'10.33"
$ws.Range("E46").Value = "  +2.04%  "
$ws.Range("E47").Value = "  +11.60%  "
$ws.Range("D48").Value = "'2.89"
$ws.Range("E48").Value = "  +1.16%  "
$ws.Range("D49").Value = "'1.52"
$ws.Range("E49").Value = "  +1.48%  "
$ws.Range("E50").Value = "  +1.94%  "
$ws.Range("D51").Value = "'72.74"
$ws.Range("E51").Value = "  +5.33%  "
